$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the shuffled weekly data (row permutation) for columns D,H,I,J,K,L,M,N,O,P,Q
$ws.Range('D2').Value = 44446
$ws.Range('H2').Value = 'Zafiro rojo'
$ws.Range('I2').Value = 'Primera'
$ws.Range('J2').Value = 10
$ws.Range('K2').Value = 34000
$ws.Range('L2').Value = 34000
$ws.Range('M2').Value = 34000
$ws.Range('N2').Value = '$/caja 18 kilos'
$ws.Range('O2').Value = 'Provincia de Limarí'
$ws.Range('P2').Value = 1889
$ws.Range('Q2').Value = 18
$ws.Range('D3').Value = 44446
$ws.Range('H3').Value = 'Zafiro verde'
$ws.Range('I3').Value = 'Primera'
$ws.Range('J3').Value = 10
$ws.Range('K3').Value = 33000
$ws.Range('L3').Value = 33000
$ws.Range('M3').Value = 33000
$ws.Range('N3').Value = '$/caja 18 kilos'
$ws.Range('O3').Value = 'Provincia de Limarí'
$ws.Range('P3').Value = 1833
$ws.Range('Q3').Value = 18
$ws.Range('D4').Value = 44446
$ws.Range('H4').Value = 'Zafiro verde'
$ws.Range('I4').Value = 'Segunda'
$ws.Range('J4').Value = 8
$ws.Range('K4').Value = 31000
$ws.Range('L4').Value = 31000
$ws.Range('M4').Value = 31000
$ws.Range('N4').Value = '$/caja 18 kilos'
$ws.Range('O4').Value = 'Provincia de Limarí'
$ws.Range('P4').Value = 1722
$ws.Range('Q4').Value = 18
$ws.Range('D5').Value = 44446
$ws.Range('H5').Value = 'Zafiro verde'
$ws.Range('I5').Value = 'Tercera'
$ws.Range('J5').Value = 12
$ws.Range('K5').Value = 29000
$ws.Range('L5').Value = 29000
$ws.Range('M5').Value = 29000
$ws.Range('N5').Value = '$/caja 18 kilos'
$ws.Range('O5').Value = 'Provincia de Limarí'
$ws.Range('P5').Value = 1611
$ws.Range('Q5').Value = 18
$ws.Range('D6').Value = 44333
$ws.Range('H6').Value = 'Zafiro verde'
$ws.Range('I6').Value = 'Primera'
$ws.Range('J6').Value = 25
$ws.Range('K6').Value = 12000
$ws.Range('L6').Value = 13000
$ws.Range('M6').Value = 12600
$ws.Range('N6').Value = '$/caja 18 kilos'
$ws.Range('O6').Value = 'Provincia de Limarí'
$ws.Range('P6').Value = 700
$ws.Range('Q6').Value = 18
$ws.Range('D7').Value = 44425
$ws.Range('H7').Value = 'Morrón rojo'
$ws.Range('I7').Value = 'Primera'
$ws.Range('J7').Value = 8
$ws.Range('K7').Value = 38000
$ws.Range('L7').Value = 38000
$ws.Range('M7').Value = 38000
$ws.Range('N7').Value = '$/caja 18 kilos'
$ws.Range('O7').Value = 'Provincia de Limarí'
$ws.Range('P7').Value = 2111
$ws.Range('Q7').Value = 18
$ws.Range('D8').Value = 44425
$ws.Range('H8').Value = 'Morrón rojo'
$ws.Range('I8').Value = 'Segunda'
$ws.Range('J8').Value = 15
$ws.Range('K8').Value = 33000
$ws.Range('L8').Value = 33000
$ws.Range('M8').Value = 33000
$ws.Range('N8').Value = '$/caja 18 kilos'
$ws.Range('O8').Value = 'Provincia de Limarí'
$ws.Range('P8').Value = 1833
$ws.Range('Q8').Value = 18
$ws.Range('D9').Value = 44425
$ws.Range('H9').Value = 'Morrón rojo'
$ws.Range('I9').Value = 'Tercera'
$ws.Range('J9').Value = 10
$ws.Range('K9').Value = 31000
$ws.Range('L9').Value = 31000
$ws.Range('M9').Value = 31000
$ws.Range('N9').Value = '$/caja 18 kilos'
$ws.Range('O9').Value = 'Provincia de Limarí'
$ws.Range('P9').Value = 1722
$ws.Range('Q9').Value = 18
$ws.Range('D10').Value = 44425
$ws.Range('H10').Value = 'Zafiro verde'
$ws.Range('I10').Value = 'Primera'
$ws.Range('J10').Value = 20
$ws.Range('K10').Value = 35000
$ws.Range('L10').Value = 35000
$ws.Range('M10').Value = 35000
$ws.Range('N10').Value = '$/caja 18 kilos'
$ws.Range('O10').Value = 'Provincia de Limarí'
$ws.Range('P10').Value = 1944
$ws.Range('Q10').Value = 18
$ws.Range('D11').Value = 44343
$ws.Range('H11').Value = 'Zafiro verde'
$ws.Range('I11').Value = 'Primera'
$ws.Range('J11').Value = 25
$ws.Range('K11').Value = 14000
$ws.Range('L11').Value = 14000
$ws.Range('M11').Value = 14000
$ws.Range('N11').Value = '$/caja 18 kilos'
$ws.Range('O11').Value = 'Provincia de Limarí'
$ws.Range('P11').Value = 778
$ws.Range('Q11').Value = 18
$ws.Range('D12').Value = 44343
$ws.Range('H12').Value = 'Zafiro verde'
$ws.Range('I12').Value = 'Segunda'
$ws.Range('J12').Value = 15
$ws.Range('K12').Value = 12000
$ws.Range('L12').Value = 12000
$ws.Range('M12').Value = 12000
$ws.Range('N12').Value = '$/caja 18 kilos'
$ws.Range('O12').Value = 'Provincia de Limarí'
$ws.Range('P12').Value = 667
$ws.Range('Q12').Value = 18
$ws.Range('D13').Value = 44421
$ws.Range('H13').Value = 'Zafiro rojo'
$ws.Range('I13').Value = 'Primera'
$ws.Range('J13').Value = 15
$ws.Range('K13').Value = 28000
$ws.Range('L13').Value = 28000
$ws.Range('M13').Value = 28000
$ws.Range('N13').Value = '$/caja 18 kilos'
$ws.Range('O13').Value = 'Provincia de Limarí'
$ws.Range('P13').Value = 1556
$ws.Range('Q13').Value = 18
$ws.Range('D14').Value = 44421
$ws.Range('H14').Value = 'Zafiro rojo'
$ws.Range('I14').Value = 'Segunda'
$ws.Range('J14').Value = 20
$ws.Range('K14').Value = 26000
$ws.Range('L14').Value = 26000
$ws.Range('M14').Value = 26000
$ws.Range('N14').Value = '$/caja 18 kilos'
$ws.Range('O14').Value = 'Provincia de Limarí'
$ws.Range('P14').Value = 1444
$ws.Range('Q14').Value = 18
$ws.Range('D15').Value = 44421
$ws.Range('H15').Value = 'Zafiro verde'
$ws.Range('I15').Value = 'Primera'
$ws.Range('J15').Value = 15
$ws.Range('K15').Value = 32000
$ws.Range('L15').Value = 32000
$ws.Range('M15').Value = 32000
$ws.Range('N15').Value = '$/caja 18 kilos'
$ws.Range('O15').Value = 'Provincia de Limarí'
$ws.Range('P15').Value = 1778
$ws.Range('Q15').Value = 18
$ws.Range('D16').Value = 44421
$ws.Range('H16').Value = 'Zafiro verde'
$ws.Range('I16').Value = 'Segunda'
$ws.Range('J16').Value = 10
$ws.Range('K16').Value = 30000
$ws.Range('L16').Value = 30000
$ws.Range('M16').Value = 30000
$ws.Range('N16').Value = '$/caja 18 kilos'
$ws.Range('O16').Value = 'Provincia de Limarí'
$ws.Range('P16').Value = 1667
$ws.Range('Q16').Value = 18
$ws.Range('D17').Value = 44421
$ws.Range('H17').Value = 'Zafiro verde'
$ws.Range('I17').Value = 'Tercera'
$ws.Range('J17').Value = 12
$ws.Range('K17').Value = 28000
$ws.Range('L17').Value = 28000
$ws.Range('M17').Value = 28000
$ws.Range('N17').Value = '$/caja 18 kilos'
$ws.Range('O17').Value = 'Provincia de Limarí'
$ws.Range('P17').Value = 1556
$ws.Range('Q17').Value = 18
$ws.Range('D18').Value = 44467
$ws.Range('H18').Value = 'Cuatro cascos verde'
$ws.Range('I18').Value = 'Primera'
$ws.Range('J18').Value = 15
$ws.Range('K18').Value = 44000
$ws.Range('L18').Value = 44000
$ws.Range('M18').Value = 44000
$ws.Range('N18').Value = '$/caja 18 kilos'
$ws.Range('O18').Value = 'Provincia de Limarí'
$ws.Range('P18').Value = 2444
$ws.Range('Q18').Value = 18
$ws.Range('D19').Value = 44467
$ws.Range('H19').Value = 'Cuatro cascos verde'
$ws.Range('I19').Value = 'Segunda'
$ws.Range('J19').Value = 25
$ws.Range('K19').Value = 42000
$ws.Range('L19').Value = 42000
$ws.Range('M19').Value = 42000
$ws.Range('N19').Value = '$/caja 18 kilos'
$ws.Range('O19').Value = 'Provincia de Limarí'
$ws.Range('P19').Value = 2333
$ws.Range('Q19').Value = 18
$ws.Range('D20').Value = 44467
$ws.Range('H20').Value = 'Morrón rojo'
$ws.Range('I20').Value = 'Primera'
$ws.Range('J20').Value = 10
$ws.Range('K20').Value = 45000
$ws.Range('L20').Value = 45000
$ws.Range('M20').Value = 45000
$ws.Range('N20').Value = '$/caja 18 kilos'
$ws.Range('O20').Value = 'Provincia de Limarí'
$ws.Range('P20').Value = 2500
$ws.Range('Q20').Value = 18
$ws.Range('D21').Value = 44467
$ws.Range('H21').Value = 'Morrón rojo'
$ws.Range('I21').Value = 'Segunda'
$ws.Range('J21').Value = 15
$ws.Range('K21').Value = 43000
$ws.Range('L21').Value = 43000
$ws.Range('M21').Value = 43000
$ws.Range('N21').Value = '$/caja 18 kilos'
$ws.Range('O21').Value = 'Provincia de Limarí'
$ws.Range('P21').Value = 2389
$ws.Range('Q21').Value = 18
$ws.Range('D22').Value = 44460
$ws.Range('H22').Value = 'Zafiro verde'
$ws.Range('I22').Value = 'Primera'
$ws.Range('J22').Value = 55
$ws.Range('K22').Value = 36000
$ws.Range('L22').Value = 36000
$ws.Range('M22').Value = 36000
$ws.Range('N22').Value = '$/caja 18 kilos'
$ws.Range('O22').Value = 'Provincia de Limarí'
$ws.Range('P22').Value = 2000
$ws.Range('Q22').Value = 18
$ws.Range('D23').Value = 44291
$ws.Range('H23').Value = 'Morrón rojo'
$ws.Range('I23').Value = 'Primera'
$ws.Range('J23').Value = 20
$ws.Range('K23').Value = 10000
$ws.Range('L23').Value = 10000
$ws.Range('M23').Value = 10000
$ws.Range('N23').Value = '$/caja 18 kilos'
$ws.Range('O23').Value = 'Provincia de Limarí'
$ws.Range('P23').Value = 556
$ws.Range('Q23').Value = 18
$ws.Range('D24').Value = 44193
$ws.Range('H24').Value = 'Zafiro rojo'
$ws.Range('I24').Value = 'Tercera'
$ws.Range('J24').Value = 20
$ws.Range('K24').Value = 28000
$ws.Range('L24').Value = 28000
$ws.Range('M24').Value = 28000
$ws.Range('N24').Value = '$/caja 18 kilos'
$ws.Range('O24').Value = 'Provincia de Limarí'
$ws.Range('P24').Value = 1556
$ws.Range('Q24').Value = 18
$ws.Range('D25').Value = 44193
$ws.Range('H25').Value = 'Zafiro verde'
$ws.Range('I25').Value = 'Primera'
$ws.Range('J25').Value = 15
$ws.Range('K25').Value = 18000
$ws.Range('L25').Value = 18000
$ws.Range('M25').Value = 18000
$ws.Range('N25').Value = '$/caja 18 kilos'
$ws.Range('O25').Value = 'Provincia de Limarí'
$ws.Range('P25').Value = 1000
$ws.Range('Q25').Value = 18
$ws.Range('D26').Value = 44193
$ws.Range('H26').Value = 'Zafiro verde'
$ws.Range('I26').Value = 'Segunda'
$ws.Range('J26').Value = 18
$ws.Range('K26').Value = 16000
$ws.Range('L26').Value = 16000
$ws.Range('M26').Value = 16000
$ws.Range('N26').Value = '$/caja 18 kilos'
$ws.Range('O26').Value = 'Provincia de Limarí'
$ws.Range('P26').Value = 889
$ws.Range('Q26').Value = 18
$ws.Range('D27').Value = 44312
$ws.Range('H27').Value = 'Zafiro rojo'
$ws.Range('I27').Value = 'Primera'
$ws.Range('J27').Value = 25
$ws.Range('K27').Value = 20000
$ws.Range('L27').Value = 20000
$ws.Range('M27').Value = 20000
$ws.Range('N27').Value = '$/caja 18 kilos'
$ws.Range('O27').Value = 'Provincia de Limarí'
$ws.Range('P27').Value = 1111
$ws.Range('Q27').Value = 18
$ws.Range('D28').Value = 44312
$ws.Range('H28').Value = 'Zafiro verde'
$ws.Range('I28').Value = 'Primera'
$ws.Range('J28').Value = 30
$ws.Range('K28').Value = 15000
$ws.Range('L28').Value = 15000
$ws.Range('M28').Value = 15000
$ws.Range('N28').Value = '$/caja 18 kilos'
$ws.Range('O28').Value = 'Provincia de Limarí'
$ws.Range('P28').Value = 833
$ws.Range('Q28').Value = 18
$ws.Range('D29').Value = 44449
$ws.Range('H29').Value = 'Zafiro rojo'
$ws.Range('I29').Value = 'Tercera'
$ws.Range('J29').Value = 25
$ws.Range('K29').Value = 34000
$ws.Range('L29').Value = 34000
$ws.Range('M29').Value = 34000
$ws.Range('N29').Value = '$/caja 18 kilos'
$ws.Range('O29').Value = 'Provincia de Limarí'
$ws.Range('P29').Value = 1889
$ws.Range('Q29').Value = 18
$ws.Range('D30').Value = 44449
$ws.Range('H30').Value = 'Zafiro verde'
$ws.Range('I30').Value = 'Primera'
$ws.Range('J30').Value = 15
$ws.Range('K30').Value = 38000
$ws.Range('L30').Value = 38000
$ws.Range('M30').Value = 38000
$ws.Range('N30').Value = '$/caja 18 kilos'
$ws.Range('O30').Value = 'Provincia de Limarí'
$ws.Range('P30').Value = 2111
$ws.Range('Q30').Value = 18
$ws.Range('D31').Value = 44449
$ws.Range('H31').Value = 'Zafiro verde'
$ws.Range('I31').Value = 'Segunda'
$ws.Range('J31').Value = 20
$ws.Range('K31').Value = 36000
$ws.Range('L31').Value = 36000
$ws.Range('M31').Value = 36000
$ws.Range('N31').Value = '$/caja 18 kilos'
$ws.Range('O31').Value = 'Provincia de Limarí'
$ws.Range('P31').Value = 2000
$ws.Range('Q31').Value = 18
$ws.Range('D32').Value = 44449
$ws.Range('H32').Value = 'Zafiro verde'
$ws.Range('I32').Value = 'Tercera'
$ws.Range('J32').Value = 25
$ws.Range('K32').Value = 34000
$ws.Range('L32').Value = 34000
$ws.Range('M32').Value = 34000
$ws.Range('N32').Value = '$/caja 18 kilos'
$ws.Range('O32').Value = 'Provincia de Limarí'
$ws.Range('P32').Value = 1889
$ws.Range('Q32').Value = 18
$ws.Range('D33').Value = 44249
$ws.Range('H33').Value = 'Cuatro cascos verde'
$ws.Range('I33').Value = 'Segunda'
$ws.Range('J33').Value = 20
$ws.Range('K33').Value = 6000
$ws.Range('L33').Value = 6000
$ws.Range('M33').Value = 6000
$ws.Range('N33').Value = '$/caja 18 kilos'
$ws.Range('O33').Value = 'Provincia de Limarí'
$ws.Range('P33').Value = 333
$ws.Range('Q33').Value = 18
$ws.Range('D34').Value = 44249
$ws.Range('H34').Value = 'Morrón rojo'
$ws.Range('I34').Value = 'Primera'
$ws.Range('J34').Value = 18
$ws.Range('K34').Value = 15000
$ws.Range('L34').Value = 15000
$ws.Range('M34').Value = 15000
$ws.Range('N34').Value = '$/caja 18 kilos'
$ws.Range('O34').Value = 'Provincia de Limarí'
$ws.Range('P34').Value = 833
$ws.Range('Q34').Value = 18
$ws.Range('D35').Value = 44186
$ws.Range('H35').Value = 'Zafiro verde'
$ws.Range('I35').Value = 'Primera'
$ws.Range('J35').Value = 20
$ws.Range('K35').Value = 17000
$ws.Range('L35').Value = 17000
$ws.Range('M35').Value = 17000
$ws.Range('N35').Value = '$/caja 18 kilos'
$ws.Range('O35').Value = 'Provincia de Limarí'
$ws.Range('P35').Value = 944
$ws.Range('Q35').Value = 18
$ws.Range('D36').Value = 44376
$ws.Range('H36').Value = 'Zafiro verde'
$ws.Range('I36').Value = 'Primera'
$ws.Range('J36').Value = 30
$ws.Range('K36').Value = 16000
$ws.Range('L36').Value = 16000
$ws.Range('M36').Value = 16000
$ws.Range('N36').Value = '$/caja 18 kilos'
$ws.Range('O36').Value = 'Provincia de Limarí'
$ws.Range('P36').Value = 889
$ws.Range('Q36').Value = 18
$ws.Range('D37').Value = 44179
$ws.Range('H37').Value = 'Zafiro rojo'
$ws.Range('I37').Value = 'Tercera'
$ws.Range('J37').Value = 20
$ws.Range('K37').Value = 22000
$ws.Range('L37').Value = 22000
$ws.Range('M37').Value = 22000
$ws.Range('N37').Value = '$/caja 18 kilos'
$ws.Range('O37').Value = 'Provincia de Limarí'
$ws.Range('P37').Value = 1222
$ws.Range('Q37').Value = 18
$ws.Range('D38').Value = 44179
$ws.Range('H38').Value = 'Zafiro verde'
$ws.Range('I38').Value = 'Primera'
$ws.Range('J38').Value = 25
$ws.Range('K38').Value = 16000
$ws.Range('L38').Value = 16000
$ws.Range('M38').Value = 16000
$ws.Range('N38').Value = '$/caja 18 kilos'
$ws.Range('O38').Value = 'Provincia de Limarí'
$ws.Range('P38').Value = 889
$ws.Range('Q38').Value = 18
$ws.Range('D39').Value = 44179
$ws.Range('H39').Value = 'Zafiro verde'
$ws.Range('I39').Value = 'Segunda'
$ws.Range('J39').Value = 20
$ws.Range('K39').Value = 14000
$ws.Range('L39').Value = 14000
$ws.Range('M39').Value = 14000
$ws.Range('N39').Value = '$/caja 18 kilos'
$ws.Range('O39').Value = 'Provincia de Limarí'
$ws.Range('P39').Value = 778
$ws.Range('Q39').Value = 18
$ws.Range('D40').Value = 44305
$ws.Range('H40').Value = 'Zafiro rojo'
$ws.Range('I40').Value = 'Primera'
$ws.Range('J40').Value = 35
$ws.Range('K40').Value = 17000
$ws.Range('L40').Value = 18000
$ws.Range('M40').Value = 17571
$ws.Range('N40').Value = '$/caja 18 kilos'
$ws.Range('O40').Value = 'Provincia de Limarí'
$ws.Range('P40').Value = 976
$ws.Range('Q40').Value = 18
$ws.Range('D41').Value = 44305
$ws.Range('H41').Value = 'Zafiro verde'
$ws.Range('I41').Value = 'Primera'
$ws.Range('J41').Value = 60
$ws.Range('K41').Value = 14000
$ws.Range('L41').Value = 14000
$ws.Range('M41').Value = 14000
$ws.Range('N41').Value = '$/caja 18 kilos'
$ws.Range('O41').Value = 'Provincia de Limarí'
$ws.Range('P41').Value = 778
$ws.Range('Q41').Value = 18
$ws.Range('D42').Value = 44232
$ws.Range('H42').Value = 'Cuatro cascos verde'
$ws.Range('I42').Value = 'Primera'
$ws.Range('J42').Value = 70
$ws.Range('K42').Value = 12000
$ws.Range('L42').Value = 12000
$ws.Range('M42').Value = 12000
$ws.Range('N42').Value = '$/caja 15 kilos'
$ws.Range('O42').Value = 'Región de Arica y Parinacota'
$ws.Range('P42').Value = 800
$ws.Range('Q42').Value = 15
$ws.Range('D43').Value = 44165
$ws.Range('H43').Value = 'Zafiro rojo'
$ws.Range('I43').Value = 'Primera'
$ws.Range('J43').Value = 125
$ws.Range('K43').Value = 30000
$ws.Range('L43').Value = 30000
$ws.Range('M43').Value = 30000
$ws.Range('N43').Value = '$/caja 15 kilos'
$ws.Range('O43').Value = 'Región de Arica y Parinacota'
$ws.Range('P43').Value = 2000
$ws.Range('Q43').Value = 15
$ws.Range('D44').Value = 44165
$ws.Range('H44').Value = 'Zafiro verde'
$ws.Range('I44').Value = 'Primera'
$ws.Range('J44').Value = 60
$ws.Range('K44').Value = 28000
$ws.Range('L44').Value = 28000
$ws.Range('M44').Value = 28000
$ws.Range('N44').Value = '$/caja 18 kilos'
$ws.Range('O44').Value = 'Provincia de Quillota'
$ws.Range('P44').Value = 1556
$ws.Range('Q44').Value = 18
$ws.Range('D45').Value = 44165
$ws.Range('H45').Value = 'Zafiro verde'
$ws.Range('I45').Value = 'Segunda'
$ws.Range('J45').Value = 40
$ws.Range('K45').Value = 25000
$ws.Range('L45').Value = 25000
$ws.Range('M45').Value = 25000
$ws.Range('N45').Value = '$/caja 18 kilos'
$ws.Range('O45').Value = 'Provincia de Quillota'
$ws.Range('P45').Value = 1389
$ws.Range('Q45').Value = 18
$ws.Range('D46').Value = 44243
$ws.Range('H46').Value = 'Cuatro cascos rojo'
$ws.Range('I46').Value = 'Primera'
$ws.Range('J46').Value = 55
$ws.Range('K46').Value = 20000
$ws.Range('L46').Value = 22000
$ws.Range('M46').Value = 21091
$ws.Range('N46').Value = '$/caja 18 kilos'
$ws.Range('O46').Value = 'Provincia de Quillota'
$ws.Range('P46').Value = 1172
$ws.Range('Q46').Value = 18
$ws.Range('D47').Value = 44243
$ws.Range('H47').Value = 'Cuatro cascos verde'
$ws.Range('I47').Value = 'Primera'
$ws.Range('J47').Value = 90
$ws.Range('K47').Value = 12000
$ws.Range('L47').Value = 13000
$ws.Range('M47').Value = 12556
$ws.Range('N47').Value = '$/caja 18 kilos'
$ws.Range('O47').Value = 'Provincia de Quillota'
$ws.Range('P47').Value = 698
$ws.Range('Q47').Value = 18
$ws.Range('D48').Value = 44236
$ws.Range('H48').Value = 'Cuatro cascos rojo'
$ws.Range('I48').Value = 'Extra'
$ws.Range('J48').Value = 60
$ws.Range('K48').Value = 25000
$ws.Range('L48').Value = 25000
$ws.Range('M48').Value = 25000
$ws.Range('N48').Value = '$/caja 18 kilos'
$ws.Range('O48').Value = 'Provincia de Limarí'
$ws.Range('P48').Value = 1389
$ws.Range('Q48').Value = 18
$ws.Range('D49').Value = 44236
$ws.Range('H49').Value = 'Cuatro cascos rojo'
$ws.Range('I49').Value = 'Primera'
$ws.Range('J49').Value = 120
$ws.Range('K49').Value = 23000
$ws.Range('L49').Value = 23000
$ws.Range('M49').Value = 23000
$ws.Range('N49').Value = '$/caja 18 kilos'
$ws.Range('O49').Value = 'Provincia de Limarí'
$ws.Range('P49').Value = 1278
$ws.Range('Q49').Value = 18
$ws.Range('D50').Value = 44236
$ws.Range('H50').Value = 'Cuatro cascos rojo'
$ws.Range('I50').Value = 'Segunda'
$ws.Range('J50').Value = 80
$ws.Range('K50').Value = 21000
$ws.Range('L50').Value = 21000
$ws.Range('M50').Value = 21000
$ws.Range('N50').Value = '$/caja 18 kilos'
$ws.Range('O50').Value = 'Provincia de Limarí'
$ws.Range('P50').Value = 1167
$ws.Range('Q50').Value = 18
$ws.Range('D51').Value = 44236
$ws.Range('H51').Value = 'Cuatro cascos rojo'
$ws.Range('I51').Value = 'Tercera'
$ws.Range('J51').Value = 50
$ws.Range('K51').Value = 18000
$ws.Range('L51').Value = 18000
$ws.Range('M51').Value = 18000
$ws.Range('N51').Value = '$/caja 18 kilos'
$ws.Range('O51').Value = 'Provincia de Limarí'
$ws.Range('P51').Value = 1000
$ws.Range('Q51').Value = 18
$ws.Range('D52').Value = 44236
$ws.Range('H52').Value = 'Cuatro cascos verde'
$ws.Range('I52').Value = 'Primera'
$ws.Range('J52').Value = 150
$ws.Range('K52').Value = 10000
$ws.Range('L52').Value = 10000
$ws.Range('M52').Value = 10000
$ws.Range('N52').Value = '$/caja 18 kilos'
$ws.Range('O52').Value = 'Provincia de Limarí'
$ws.Range('P52').Value = 556
$ws.Range('Q52').Value = 18
$ws.Range('D53').Value = 44236
$ws.Range('H53').Value = 'Cuatro cascos verde'
$ws.Range('I53').Value = 'Segunda'
$ws.Range('J53').Value = 100
$ws.Range('K53').Value = 8000
$ws.Range('L53').Value = 8000
$ws.Range('M53').Value = 8000
$ws.Range('N53').Value = '$/caja 18 kilos'
$ws.Range('O53').Value = 'Provincia de Limarí'
$ws.Range('P53').Value = 444
$ws.Range('Q53').Value = 18
$ws.Range('D54').Value = 44236
$ws.Range('H54').Value = 'Cuatro cascos verde'
$ws.Range('I54').Value = 'Tercera'
$ws.Range('J54').Value = 75
$ws.Range('K54').Value = 6000
$ws.Range('L54').Value = 6000
$ws.Range('M54').Value = 6000
$ws.Range('N54').Value = '$/caja 18 kilos'
$ws.Range('O54').Value = 'Provincia de Limarí'
$ws.Range('P54').Value = 333
$ws.Range('Q54').Value = 18
$ws.Range('D55').Value = 44277
$ws.Range('H55').Value = 'Zafiro rojo'
$ws.Range('I55').Value = 'Primera'
$ws.Range('J55').Value = 20
$ws.Range('K55').Value = 16000
$ws.Range('L55').Value = 16000
$ws.Range('M55').Value = 16000
$ws.Range('N55').Value = '$/caja 18 kilos'
$ws.Range('O55').Value = 'Provincia de Limarí'
$ws.Range('P55').Value = 889
$ws.Range('Q55').Value = 18
$ws.Range('D56').Value = 44277
$ws.Range('H56').Value = 'Zafiro verde'
$ws.Range('I56').Value = 'Primera'
$ws.Range('J56').Value = 25
$ws.Range('K56').Value = 12000
$ws.Range('L56').Value = 12000
$ws.Range('M56').Value = 12000
$ws.Range('N56').Value = '$/caja 18 kilos'
$ws.Range('O56').Value = 'Provincia de Limarí'
$ws.Range('P56').Value = 667
$ws.Range('Q56').Value = 18
$ws.Range('D57').Value = 44432
$ws.Range('H57').Value = 'Zafiro verde'
$ws.Range('I57').Value = 'Primera'
$ws.Range('J57').Value = 20
$ws.Range('K57').Value = 35000
$ws.Range('L57').Value = 35000
$ws.Range('M57').Value = 35000
$ws.Range('N57').Value = '$/caja 18 kilos'
$ws.Range('O57').Value = 'Provincia de Limarí'
$ws.Range('P57').Value = 1944
$ws.Range('Q57').Value = 18
$ws.Range('D58').Value = 44284
$ws.Range('H58').Value = 'Zafiro rojo'
$ws.Range('I58').Value = 'Primera'
$ws.Range('J58').Value = 20
$ws.Range('K58').Value = 16000
$ws.Range('L58').Value = 16000
$ws.Range('M58').Value = 16000
$ws.Range('N58').Value = '$/caja 18 kilos'
$ws.Range('O58').Value = 'Provincia de Limarí'
$ws.Range('P58').Value = 889
$ws.Range('Q58').Value = 18
$ws.Range('D59').Value = 44284
$ws.Range('H59').Value = 'Zafiro verde'
$ws.Range('I59').Value = 'Primera'
$ws.Range('J59').Value = 25
$ws.Range('K59').Value = 12000
$ws.Range('L59').Value = 12000
$ws.Range('M59').Value = 12000
$ws.Range('N59').Value = '$/caja 18 kilos'
$ws.Range('O59').Value = 'Provincia de Limarí'
$ws.Range('P59').Value = 667
$ws.Range('Q59').Value = 18
$ws.Range('D60').Value = 44435
$ws.Range('H60').Value = 'Zafiro verde'
$ws.Range('I60').Value = 'Primera'
$ws.Range('J60').Value = 30
$ws.Range('K60').Value = 35000
$ws.Range('L60').Value = 35000
$ws.Range('M60').Value = 35000
$ws.Range('N60').Value = '$/caja 18 kilos'
$ws.Range('O60').Value = 'Provincia de Limarí'
$ws.Range('P60').Value = 1944
$ws.Range('Q60').Value = 18
$ws.Range('D61').Value = 44200
$ws.Range('H61').Value = 'Zafiro verde'
$ws.Range('I61').Value = 'Primera'
$ws.Range('J61').Value = 15
$ws.Range('K61').Value = 16000
$ws.Range('L61').Value = 16000
$ws.Range('M61').Value = 16000
$ws.Range('N61').Value = '$/caja 18 kilos'
$ws.Range('O61').Value = 'Provincia de Limarí'
$ws.Range('P61').Value = 889
$ws.Range('Q61').Value = 18
$ws.Range('D62').Value = 44200
$ws.Range('H62').Value = 'Zafiro verde'
$ws.Range('I62').Value = 'Segunda'
$ws.Range('J62').Value = 10
$ws.Range('K62').Value = 13000
$ws.Range('L62').Value = 13000
$ws.Range('M62').Value = 13000
$ws.Range('N62').Value = '$/caja 18 kilos'
$ws.Range('O62').Value = 'Provincia de Limarí'
$ws.Range('P62').Value = 722
$ws.Range('Q62').Value = 18
$ws.Range('D63').Value = 44280
$ws.Range('H63').Value = 'Zafiro rojo'
$ws.Range('I63').Value = 'Primera'
$ws.Range('J63').Value = 30
$ws.Range('K63').Value = 16000
$ws.Range('L63').Value = 16000
$ws.Range('M63').Value = 16000
$ws.Range('N63').Value = '$/caja 18 kilos'
$ws.Range('O63').Value = 'Provincia de Limarí'
$ws.Range('P63').Value = 889
$ws.Range('Q63').Value = 18
$ws.Range('D64').Value = 44280
$ws.Range('H64').Value = 'Zafiro verde'
$ws.Range('I64').Value = 'Primera'
$ws.Range('J64').Value = 45
$ws.Range('K64').Value = 12000
$ws.Range('L64').Value = 12000
$ws.Range('M64').Value = 12000
$ws.Range('N64').Value = '$/caja 18 kilos'
$ws.Range('O64').Value = 'Provincia de Limarí'
$ws.Range('P64').Value = 667
$ws.Range('Q64').Value = 18
$ws.Range('D65').Value = 44418
$ws.Range('H65').Value = 'Morrón rojo'
$ws.Range('I65').Value = 'Primera'
$ws.Range('J65').Value = 10
$ws.Range('K65').Value = 33000
$ws.Range('L65').Value = 33000
$ws.Range('M65').Value = 33000
$ws.Range('N65').Value = '$/caja 18 kilos'
$ws.Range('O65').Value = 'Provincia de Limarí'
$ws.Range('P65').Value = 1833
$ws.Range('Q65').Value = 18
$ws.Range('D66').Value = 44418
$ws.Range('H66').Value = 'Morrón rojo'
$ws.Range('I66').Value = 'Segunda'
$ws.Range('J66').Value = 12
$ws.Range('K66').Value = 31000
$ws.Range('L66').Value = 31000
$ws.Range('M66').Value = 31000
$ws.Range('N66').Value = '$/caja 18 kilos'
$ws.Range('O66').Value = 'Provincia de Limarí'
$ws.Range('P66').Value = 1722
$ws.Range('Q66').Value = 18
$ws.Range('D67').Value = 44418
$ws.Range('H67').Value = 'Zafiro rojo'
$ws.Range('I67').Value = 'Primera'
$ws.Range('J67').Value = 10
$ws.Range('K67').Value = 28000
$ws.Range('L67').Value = 28000
$ws.Range('M67').Value = 28000
$ws.Range('N67').Value = '$/caja 18 kilos'
$ws.Range('O67').Value = 'Provincia de Limarí'
$ws.Range('P67').Value = 1556
$ws.Range('Q67').Value = 18
$ws.Range('D68').Value = 44418
$ws.Range('H68').Value = 'Zafiro rojo'
$ws.Range('I68').Value = 'Segunda'
$ws.Range('J68').Value = 15
$ws.Range('K68').Value = 26000
$ws.Range('L68').Value = 26000
$ws.Range('M68').Value = 26000
$ws.Range('N68').Value = '$/caja 18 kilos'
$ws.Range('O68').Value = 'Provincia de Limarí'
$ws.Range('P68').Value = 1444
$ws.Range('Q68').Value = 18
$ws.Range('D69').Value = 44418
$ws.Range('H69').Value = 'Zafiro verde'
$ws.Range('I69').Value = 'Primera'
$ws.Range('J69').Value = 10
$ws.Range('K69').Value = 28000
$ws.Range('L69').Value = 28000
$ws.Range('M69').Value = 28000
$ws.Range('N69').Value = '$/caja 18 kilos'
$ws.Range('O69').Value = 'Provincia de Limarí'
$ws.Range('P69').Value = 1556
$ws.Range('Q69').Value = 18
$ws.Range('D70').Value = 44270
$ws.Range('H70').Value = 'Zafiro rojo'
$ws.Range('I70').Value = 'Primera'
$ws.Range('J70').Value = 25
$ws.Range('K70').Value = 17000
$ws.Range('L70').Value = 17000
$ws.Range('M70').Value = 17000
$ws.Range('N70').Value = '$/caja 18 kilos'
$ws.Range('O70').Value = 'Provincia de Limarí'
$ws.Range('P70').Value = 944
$ws.Range('Q70').Value = 18
$ws.Range('D71').Value = 44270
$ws.Range('H71').Value = 'Zafiro verde'
$ws.Range('I71').Value = 'Primera'
$ws.Range('J71').Value = 30
$ws.Range('K71').Value = 10000
$ws.Range('L71').Value = 10000
$ws.Range('M71').Value = 10000
$ws.Range('N71').Value = '$/caja 18 kilos'
$ws.Range('O71').Value = 'Provincia de Limarí'
$ws.Range('P71').Value = 556
$ws.Range('Q71').Value = 18
$ws.Range('D72').Value = 44273
$ws.Range('H72').Value = 'Zafiro rojo'
$ws.Range('I72').Value = 'Primera'
$ws.Range('J72').Value = 45
$ws.Range('K72').Value = 17000
$ws.Range('L72').Value = 17000
$ws.Range('M72').Value = 17000
$ws.Range('N72').Value = '$/caja 18 kilos'
$ws.Range('O72').Value = 'Provincia de Limarí'
$ws.Range('P72').Value = 944
$ws.Range('Q72').Value = 18
$ws.Range('D73').Value = 44273
$ws.Range('H73').Value = 'Zafiro rojo'
$ws.Range('I73').Value = 'Segunda'
$ws.Range('J73').Value = 30
$ws.Range('K73').Value = 15000
$ws.Range('L73').Value = 15000
$ws.Range('M73').Value = 15000
$ws.Range('N73').Value = '$/caja 18 kilos'
$ws.Range('O73').Value = 'Provincia de Limarí'
$ws.Range('P73').Value = 833
$ws.Range('Q73').Value = 18
$ws.Range('D74').Value = 44273
$ws.Range('H74').Value = 'Zafiro rojo'
$ws.Range('I74').Value = 'Tercera'
$ws.Range('J74').Value = 25
$ws.Range('K74').Value = 13000
$ws.Range('L74').Value = 13000
$ws.Range('M74').Value = 13000
$ws.Range('N74').Value = '$/caja 18 kilos'
$ws.Range('O74').Value = 'Provincia de Limarí'
$ws.Range('P74').Value = 722
$ws.Range('Q74').Value = 18
$ws.Range('D75').Value = 44350
$ws.Range('H75').Value = 'Zafiro rojo'
$ws.Range('I75').Value = 'Primera'
$ws.Range('J75').Value = 15
$ws.Range('K75').Value = 28000
$ws.Range('L75').Value = 28000
$ws.Range('M75').Value = 28000
$ws.Range('N75').Value = '$/caja 15 kilos'
$ws.Range('O75').Value = 'Región de Arica y Parinacota'
$ws.Range('P75').Value = 1867
$ws.Range('Q75').Value = 15
$ws.Range('D76').Value = 44350
$ws.Range('H76').Value = 'Zafiro verde'
$ws.Range('I76').Value = 'Primera'
$ws.Range('J76').Value = 35
$ws.Range('K76').Value = 14000
$ws.Range('L76').Value = 14000
$ws.Range('M76').Value = 14000
$ws.Range('N76').Value = '$/caja 15 kilos'
$ws.Range('O76').Value = 'Región de Arica y Parinacota'
$ws.Range('P76').Value = 933
$ws.Range('Q76').Value = 15
$ws.Range('D77').Value = 44474
$ws.Range('H77').Value = 'Cuatro cascos verde'
$ws.Range('I77').Value = 'Primera'
$ws.Range('J77').Value = 15
$ws.Range('K77').Value = 38000
$ws.Range('L77').Value = 38000
$ws.Range('M77').Value = 38000
$ws.Range('N77').Value = '$/caja 18 kilos'
$ws.Range('O77').Value = 'Provincia de Limarí'
$ws.Range('P77').Value = 2111
$ws.Range('Q77').Value = 18
$ws.Range('D78').Value = 44474
$ws.Range('H78').Value = 'Cuatro cascos verde'
$ws.Range('I78').Value = 'Segunda'
$ws.Range('J78').Value = 25
$ws.Range('K78').Value = 36000
$ws.Range('L78').Value = 36000
$ws.Range('M78').Value = 36000
$ws.Range('N78').Value = '$/caja 18 kilos'
$ws.Range('O78').Value = 'Provincia de Limarí'
$ws.Range('P78').Value = 2000
$ws.Range('Q78').Value = 18
$ws.Range('D79').Value = 44474
$ws.Range('H79').Value = 'Cuatro cascos verde'
$ws.Range('I79').Value = 'Tercera'
$ws.Range('J79').Value = 20
$ws.Range('K79').Value = 34000
$ws.Range('L79').Value = 34000
$ws.Range('M79').Value = 34000
$ws.Range('N79').Value = '$/caja 18 kilos'
$ws.Range('O79').Value = 'Provincia de Limarí'
$ws.Range('P79').Value = 1889
$ws.Range('Q79').Value = 18
